# Update the "last modified" date field on the slide master from
# 05.06.2024 to 06.06.2024 (shape "Rectangle 6" on SlideMaster1).
$p = $ppt.ActivePresentation

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "05.06.2024") {
            $shape.TextFrame.TextRange.Text = "06.06.2024"
        }
    }
}

# Mark the last five slides (32-36) as hidden in the slide show.
for ($i = 32; $i -le 36; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.SlideShowTransition.Hidden = $true
}
